$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new trade record on row 4, mirroring the existing rows' layout.
$ws.Range("A4").Value = 42633.678611111114
$ws.Range("A4").NumberFormat = "m/d/yy h:mm"

$ws.Range("B4").Value = $false

$ws.Range("C4").Value = 9980.5
$ws.Range("D4").Value = 10000
$ws.Range("E4").Value = 108.67
$ws.Range("F4").Value = 108.25

$ws.Range("G4").Value = $false
$ws.Range("G4").NumberFormat = "m/d/yy h:mm"

$ws.Range("H4").Value = -0.39

$ws.Range("I4").Value = $false
